# Apply updated allocation rule summary values (newest airtoxics NATA data)
# to the "Means" and "Standard Deviations" worksheets.

$wb = $excel.ActiveWorkbook

# --- Sheet "Means" ---
$wsMeans = $wb.Worksheets.Item("Means")

# Row 9: Total Cancer Risk (per million)
$wsMeans.Range("B9").Value = 23
$wsMeans.Range("C9").Value = 23
$wsMeans.Range("D9").Value = 30
$wsMeans.Range("E9").Value = 30
$wsMeans.Range("F9").Value = 30
$wsMeans.Range("G9").Value = 28

# Row 10: Total Respiratory (hazard quotient)
$wsMeans.Range("B10").Value = 0.27
$wsMeans.Range("C10").Value = 0.3
$wsMeans.Range("D10").Value = 0.36
$wsMeans.Range("E10").Value = 0.34
$wsMeans.Range("F10").Value = 0.34
$wsMeans.Range("G10").Value = 0.31

# --- Sheet "Standard Deviations" ---
$wsSD = $wb.Worksheets.Item("Standard Deviations")

# Row 9: Total Cancer Risk (per million)
$wsSD.Range("B9").Value = 7.2
$wsSD.Range("C9").Value = 4.7
$wsSD.Range("D9").Value = 0
$wsSD.Range("E9").Value = 0
$wsSD.Range("F9").Value = 0
$wsSD.Range("G9").Value = 4.5

# Row 10: Total Respiratory (hazard quotient)
$wsSD.Range("B10").Value = 0.094
$wsSD.Range("C10").Value = 0.022
$wsSD.Range("D10").Value = 0.058
$wsSD.Range("E10").Value = 0.053
$wsSD.Range("F10").Value = 0.05
$wsSD.Range("G10").Value = 0.035

$wb.Save()
